$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. D7: Status changes from "Resolved" to "Won't Fix" -----------------
$ws.Range("D7").Value = "Won't Fix"

# --- 2. Fill previously-blank rows 12-16 with new bug reports --------------
# First, stamp the row formatting by copying from existing fully-formatted
# rows so the cell styles (s="2"/"3"/"4"/"5") match what Excel would apply
# when a user fills in a templated blank row.
$ws.Range("A11:I11").Copy($ws.Range("A12:I12"))
$ws.Range("A11:I11").Copy($ws.Range("A13:I13"))
$ws.Range("A11:I11").Copy($ws.Range("A14:I14"))
$ws.Range("A11:I11").Copy($ws.Range("A15:I15"))
$ws.Range("A11:I11").Copy($ws.Range("A16:I16"))
# Row 15's G/H columns use the wrap-text style (like row 10's G column).
$ws.Range("G10").Copy($ws.Range("G15"))
$ws.Range("G10").Copy($ws.Range("H15"))

# Row 12
$ws.Range("A12").Value = 45388
$ws.Range("B12").Value = "Annie"
$ws.Range("C12").Value = "P3 - Annoying, could ship but loss of marks"
$ws.Range("D12").Value = "Resolved"
$ws.Range("E12").Value = "Annie"
$ws.Range("F12").Value = "Bullets and powerups display overtop of cutscenes"
$ws.Range("G12").Value = "nothing should render over top of  the cutscenes"
$ws.Range("H12").Value = "bullets and powerups over cutscenes"
$ws.Range("I12").Value = "shoot bullets as the cutscenes are begining"

# Row 13
$ws.Range("A13").Value = 45389
$ws.Range("B13").Value = "Andrew"
$ws.Range("C13").Value = "P4 - Annoying/Limiting but could ship as is"
$ws.Range("D13").Value = "Won't Fix"
$ws.Range("E13").Value = "Andrew"
$ws.Range("F13").Value = "player will keep recieveing the key borad input (super rare case)even no been pressed"
$ws.Range("G13").Value = "Key input recieve normally"
$ws.Range("H13").Value = "Key input will keep input(player keep moving in a direction)"
$ws.Range("I13").Value = "Not sure, but when it occur, just press the key in the same direcation and it will fix this issue"

# Row 14
$ws.Range("A14").Value = 45390
$ws.Range("B14").Value = "Andrew"
$ws.Range("C14").Value = "P4 - Annoying/Limiting but could ship as is"
$ws.Range("D14").Value = "Won't Fix"
$ws.Range("E14").Value = "Eddie"
$ws.Range("F14").Value = "Boss is using AABB collision "
$ws.Range("G14").Value = "Boss is using mesh"
$ws.Range("H14").Value = "Boss use AABB"
$ws.Range("I14").Value = "Shoot to  boss"

# Row 15
$ws.Range("A15").Value = 45390
$ws.Range("B15").Value = "Andrew"
$ws.Range("C15").Value = "P4 - Annoying/Limiting but could ship as is"
$ws.Range("D15").Value = "Won't Fix"
$ws.Range("E15").Value = "Andrew"
$ws.Range("F15").Value = "Split force to strong when the enemy get too cloose"
$ws.Range("G15").Value = "Split force should be mild"
$ws.Range("H15").Value = "Split force to strong when the enemy get too clooseo strong when the enemy get too cloose"
$ws.Range("I15").Value = "Try to make two enemy really close"

# Row 16
$ws.Range("A16").Value = 45390
$ws.Range("B16").Value = "Andrew"
$ws.Range("C16").Value = "P4 - Annoying/Limiting but could ship as is"
$ws.Range("D16").Value = "Won't Fix"
$ws.Range("E16").Value = "Andrew"
$ws.Range("F16").Value = "Some super minior layering issue, when the enemy get too cloose"
$ws.Range("G16").Value = "No layering issue."
$ws.Range("H16").Value = "Some layering issue"
$ws.Range("I16").Value = "try to get enemy close"
